$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.06%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.58%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.770"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.70%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08336"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.05%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.814"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.98%"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.74%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.961"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.20%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.922"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.31%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9330"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.46%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1237"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.28%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1955"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.10%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09486"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.07%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "5.51%"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.77%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001304"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.74%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005930"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-5.32%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.504"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.93%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.022"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "9.36%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.54%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2570"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.16%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04406"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.37%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.34%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004394"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.38%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001190"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.48%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003990"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.08%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02838"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.26%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05653"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.62%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007901"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.85%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.82%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009068"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.44%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002100"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.75%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009935"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.71%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007286"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.51%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.36%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003970"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10.84%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002278"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.42%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.36%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.36%"
